# edit.ps1 -- apply "update pre slides and add paper" changes
#
#  1. Bump the cached "datetimeFigureOut" field text from 2021/11/21 to
#     2021/11/22 on the slide master and on every slide layout's date
#     placeholder.
#  2. Fill in the (previously empty) subtitle placeholder on the title
#     slide (slide 1) with the author list and the presentation date.
#  3. Re-order the author names on the "Team" slide (slide 10) from
#     "<Family> <Given>" to "<Given> <Family>".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholders: slide master + every custom layout.
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -eq -1) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "2021/11/21") {
                $tr.Text = "2021/11/22"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $lay = $layouts.Item($li)
    Update-DatePlaceholder $lay.Shapes
}

# ---------------------------------------------------------------------
# 2) Title slide (slide 1): subtitle placeholder gets the author line
#    and the date line.
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$subtitleShape = $s1.Shapes.Item(2)
$subtitleRange = $subtitleShape.TextFrame.TextRange

# Build the second paragraph first (while it is the only run in the
# frame) so LanguageID correctly lands on it, then prepend the first
# paragraph -- the newly inserted text inherits the language of the
# run it now sits next to.
$subtitleRange.InsertBefore("November 22, 2021") | Out-Null
$subtitleRange.LanguageID = "en-US"

$fullRange = $subtitleShape.TextFrame.TextRange
$fullRange.InsertBefore("Yuhan Zhou  Wenrui Liu  Xiaolong Huang" + [char]13) | Out-Null

# ---------------------------------------------------------------------
# 3) Team slide (slide 10): swap "Family Given" -> "Given Family".
# ---------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$teamShape = $s10.Shapes.Item(2)
$teamRange = $teamShape.TextFrame.TextRange

$teamRange.Replace("Zhou Yuhan ", "Yuhan Zhou ") | Out-Null
$teamRange.Replace("Liu Wenrui ", "Wenrui Liu ") | Out-Null
$teamRange.Replace("Huang Xiaolong ", "Xiaolong Huang ") | Out-Null

Write-Host "edit.ps1 completed"
